# Fruta / hortaliza, semanal
# Insert a new data row at row 32 (pushing the existing rows 32:61 down to 33:62)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32 and below down by one row, keeping all existing data/formatting intact.
$ws.Rows.Item(32).Insert()

# Fill in the new record for row 32.
$ws.Cells.Item(32, 1).Value  = 1
$ws.Cells.Item(32, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value  = 44447
$ws.Cells.Item(32, 5).Value  = 15
$ws.Cells.Item(32, 6).Value  = "Fruta"
$ws.Cells.Item(32, 7).Value  = 100102
$ws.Cells.Item(32, 8).Value  = "Cítricos"
$ws.Cells.Item(32, 9).Value  = 100102005
$ws.Cells.Item(32, 10).Value = "Naranja"
$ws.Cells.Item(32, 11).Value = "Lane Late"
$ws.Cells.Item(32, 12).Value = "Segunda"
$ws.Cells.Item(32, 13).Value = 250
$ws.Cells.Item(32, 14).Value = 600
$ws.Cells.Item(32, 15).Value = 650
$ws.Cells.Item(32, 16).Value = 625
$ws.Cells.Item(32, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(32, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(32, 19).Value = 625
$ws.Cells.Item(32, 20).Value = 1
